# Trade #73 closed at 2026-02-17 08:57:51 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# of the live trading results workbook to reflect the newly closed trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.57   # Current Capital
$summary.Range("B4").Value = 0.58      # Total P&L $
$summary.Range("B6").Value = 73        # Total Trades
$summary.Range("B7").Value = 31        # Winning Trades
$summary.Range("B9").Value = 42.47     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.57     # Capital
$status.Range("D4").Value = 73         # Trades
$status.Range("E4").Value = 0.58       # P&L $
$status.Range("F4").Value = 0.57       # P&L %
$status.Range("G4").Value = 42.47      # Win Rate %

# ---------------------------------------------------------------------
# Helper that appends the new trade #73 row to a trades-log sheet.
# ---------------------------------------------------------------------
function Add-Trade73Row($sheet) {
    # Force the Date column to be stored as literal text (matches the
    # rest of the column) instead of being auto-parsed into a date serial.
    $sheet.Range("B74").NumberFormat = "@"

    $sheet.Range("A74").Value = 73
    $sheet.Range("B74").Value = "2026-02-17"
    $sheet.Range("C74").Value = "08:57:44"
    $sheet.Range("D74").Value = "MarketMaking"
    $sheet.Range("E74").Value = "UP"
    $sheet.Range("F74").Value = 0.1
    $sheet.Range("G74").Value = 0.11
    $sheet.Range("H74").Value = "CLOSED"
    $sheet.Range("I74").Value = 10
    $sheet.Range("J74").Value = 0.01
    $sheet.Range("K74").Value = 100.57
    $sheet.Range("L74").Value = 0
    $sheet.Range("M74").Value = 0
    $sheet.Range("N74").Value = 0.6
    $sheet.Range("O74").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P74").Value = "early_exit"
    $sheet.Range("Q74").Value = 0.13
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade73Row $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade73Row $marketMaking
